# Add new "Explicacion" (column L) text for App-Control-only detections.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = 'El control de aplicaciones sólo detecta el uso de SSH, siendo insuficiente'
$ws.Range("L3").Value = 'El control de aplicaciones no detecta el uso del protocolo falso TLS'
$ws.Range("L6").Value = 'Las alertas de control de aplicaciones sólo muestran los protocolos usados, pero no indagan en la ejecución del script de duplicado'
$ws.Range("L8").Value = 'El control de aplicaciones sólo indica el uso de FTP'
$ws.Range("L9").Value = 'El control de aplicaciones sólo indica el uso de SMB o funciones de compartición de archivos. No es suficiente para determinar que se está recogiendo información de usuarios del sistema'
$ws.Range("L17").Value = 'El control de aplicaciones sólo detecta el uso de Ping'
$ws.Range("L18").Value = 'El control de aplicaciones sólo detecta el uso de FTP'
$ws.Range("L29").Value = 'El control de aplicaciones no marca el uso de la funcionalidad de registro remoto'
$ws.Range("L32").Value = 'El control de aplicaciones no marca el acceso a la VPN'
$ws.Range("L33").Value = 'El control de aplicaciones no marca el uso de una extensión, sino el uso de un navegador. No se puede determinar como válida'
$ws.Range("L34").Value = 'El control de aplicaciones sólo marca la funcionalidad de compartición de archivos, nada más'
$ws.Range("L35").Value = 'El control de aplicaciones sólo marca la funcionalidad de compartición de archivos, nada más'
$ws.Range("L36").Value = 'El control de aplicaciones sólo marca el uso del navegador y las peticiones HTTP'
$ws.Range("L39").Value = 'El control de aplicaciones sólo marca el uso de RDP, nada más'
$ws.Range("L42").Value = 'El control de aplicaciones sólo marca el uso de SSH'
$ws.Range("L43").Value = 'El control de aplicaciones sólo marca el uso deL navegador Firefox'
$ws.Range("L44").Value = 'El control de aplicaciones sólo marca el uso deL navegador Firefox'
$ws.Range("L45").Value = 'El control de aplicaciones sólo marca el uso de Github, pero no es suficiente para detectar el ataque'
$ws.Range("L46").Value = 'El control de aplicaciones sólo detecta navegadores y otros'
$ws.Range("L47").Value = 'El control de aplicaciones detecta el uso de Kerberos y LDAP, pero no indica nada de un acceso a los dominios de confianza'

# Restore the author's last cursor/selection position (and scroll the
# viewport so row 45 is at the top, matching topLeftCell="A45").
$excel.ActiveWindow.ScrollRow = 45
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B73").Select()
